$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "ontology_" -> "GamingOntology_"
# ------------------------------------------------------------------
$d.Content.Find.Execute("ontology_", $false, $false, $false, $false, $false,
                         $true, 1, $false, "GamingOntology_", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Split the run "done,owl" into two runs: "done" and ",owl"
#    (no textual change, just a run split - achieved by toggling a
#    character formatting property on/off over the first half so the
#    host is forced to break the run there)
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("done,owl", $false, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng.Start + 4
$r1 = $d.Range($rng.Start, $splitPos)
$r1.Font.Bold = 1
$r1.Font.Bold = 0

# ------------------------------------------------------------------
# 3. Split the run ". A felület bekéri a keresett " into ". A f" and
#    "elület bekéri a keresett ", moving the _GoBack bookmark from the
#    end of the paragraph to right between these two new runs.
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(". A f", $false, $false, $false, $false, $false,
                    $true, 1, $false, "", 0) | Out-Null
$splitPos2 = $rng2.End

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$newBmRange = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("_GoBack", $newBmRange) | Out-Null
